# Update F-column (views/浏览) counts on three worksheets to reflect the
# output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 2747
$ws1.Range("F10").Value = 6229
$ws1.Range("F19").Value = 1510
$ws1.Range("F24").Value = 1049
$ws1.Range("F25").Value = 230
$ws1.Range("F36").Value = 1487
$ws1.Range("F38").Value = 1039
$ws1.Range("F42").Value = 283
$ws1.Range("F46").Value = 120

# --- Sheet "本地生活" --------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F8").Value  = 1468
$ws3.Range("F11").Value = 827
$ws3.Range("F12").Value = 714

# --- Sheet "全部类型" --------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value  = 2747
$ws4.Range("F10").Value = 1468
$ws4.Range("F13").Value = 6229
$ws4.Range("F14").Value = 827
$ws4.Range("F19").Value = 1510
$ws4.Range("F24").Value = 1049
$ws4.Range("F25").Value = 230
$ws4.Range("F35").Value = 1487
$ws4.Range("F36").Value = 1039
$ws4.Range("F40").Value = 283
$ws4.Range("F44").Value = 120
